$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "GCU0D88"
$ws.Range("B2").Value = "SILVIO CESAR VIERIA"
$ws.Range("C2").Value = 53.87
$ws.Range("D2").Value = "[['CAMPINAS', ['1682409', '1683810']]]"
